$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D:K data to E:L), mirroring the
# addition of a new fiscal-year (2018-12-31) column in front of the existing
# Income Statement / Balance Sheet / Cash Flow tables.
$ws.Columns("D").Insert()

# The freshly inserted column has no formatting; copy the number formats
# (and therefore the correct shared cell styles) from column E, which holds
# the data that used to live in column D before the shift.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Populate the new column (D) with the new period's figures -----------
# Period Ending header rows (dates)
$ws.Range("D7").Value = 43465
$ws.Range("D38").Value = 43465
$ws.Range("D80").Value = 43465

# Income Statement
$ws.Range("D8").Value = 2826600
$ws.Range("D9").Value = 1369700
$ws.Range("D10").Value = 1456900
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 265100
$ws.Range("D17").Value = 1632000
$ws.Range("D18").Value = 1194600
$ws.Range("D20").Value = 340000
$ws.Range("D21").Value = 1799600
$ws.Range("D22").Value = 200600
$ws.Range("D23").Value = 1334000
$ws.Range("D24").Value = 100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1333900
$ws.Range("D27").Value = 1333900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -340000
$ws.Range("D33").Value = 1333900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1333900

# Balance Sheet
$ws.Range("D41").Value = 218300
$ws.Range("E41").Value = 160800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 129200
$ws.Range("D44").Value = 185700
$ws.Range("D45").Value = 113200
$ws.Range("D46").Value = 646300
$ws.Range("E46").Value = 596600
$ws.Range("D47").Value = 1097200
$ws.Range("D48").Value = 5798200
$ws.Range("D49").Value = 104400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 101400
$ws.Range("E52").Value = 28700
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 7747500
$ws.Range("D57").Value = 138700
$ws.Range("D58").Value = 59500
$ws.Range("D59").Value = 478300
$ws.Range("D60").Value = 676600
$ws.Range("D61").Value = 4211400
$ws.Range("D62").Value = 216200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 5104100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2643400
$ws.Range("D77").Value = 0

# Cash Flow Statement
$ws.Range("D81").Value = 1333900
$ws.Range("D83").Value = 265100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1353000
$ws.Range("E89").Value = 1131200
$ws.Range("F89").Value = 973300
$ws.Range("D91").Value = -552300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -119300
$ws.Range("E94").Value = -593200
$ws.Range("F94").Value = -866600
$ws.Range("D96").Value = -865400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -1100500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 133200

Write-Host "Applied MMP financial update"
